$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for all existing data rows
#    (rows 2-131) from 2023-09-11 (45180) to 2023-09-12 (45181).
$ws.Range("C2:C131").Value = 45181

# 2. Row 131 picks up an explicit row height (15pt, custom height flag) in the
#    target file even though it matches the sheet default.
$ws.Rows.Item(131).RowHeight = 15

# 3. Append a new record as row 132.
$newRow = 132

$ws.Cells.Item($newRow, 1).Value = "A 42350-2023"

$ws.Cells.Item($newRow, 2).Value = 45180
$ws.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 3).Value = 45181
$ws.Cells.Item($newRow, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 4).Value = "VÄRMLANDS LÄN"
$ws.Cells.Item($newRow, 5).Value = "STORFORS"

$ws.Cells.Item($newRow, 7).Value = 1.6
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
$ws.Cells.Item($newRow, 11).Value = 0
$ws.Cells.Item($newRow, 12).Value = 0
$ws.Cells.Item($newRow, 13).Value = 0
$ws.Cells.Item($newRow, 14).Value = 0
$ws.Cells.Item($newRow, 15).Value = 0
$ws.Cells.Item($newRow, 16).Value = 0
$ws.Cells.Item($newRow, 17).Value = 0

$ws.Cells.Item($newRow, 18).Value = ""
$ws.Cells.Item($newRow, 18).WrapText = $true
